# Append the new "Statistical Analysis" example write-up paragraphs.
# The new content goes right after the existing "Statistical Analysis"
# heading paragraph (the last paragraph in the document body) and before
# the trailing section properties.
$d = $word.ActiveDocument

# Collapsed range positioned immediately before the final paragraph mark
# (i.e. right before <w:sectPr>), so inserting here appends new paragraphs
# after the "Statistical Analysis" heading without disturbing it.
$insertPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">Example write up (must be rewritten)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/></w:pPr><w:r><w:t xml:space="preserve">(Reporting rows from final model)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/></w:pPr><w:r><w:t xml:space="preserve">Figure 1 plots the model estimates from the GCA, and the full model summary is</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">available in Appendices 1 and 2.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">We report the results for the M group and then provide comparisons with and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">between the learner groups.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">The model intercept estimates the log odds of M fixating on the target,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">averaging over the time course, lexical stress, and syllable structure, at the</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">mean working memory (XXX).</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">The log odds were</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">γ</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t xml:space="preserve">00</w:t></w:r><w:r><w:t xml:space="preserve"> = 1.17</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">(proportion: .76).</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">The linear, quadratic, and cubic polynomial time terms captured the sigmoid</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">shape of the time course and were retained in the model</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">(</w:t></w:r><w:r><w:t xml:space="preserve">γ</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t xml:space="preserve">10</w:t></w:r><w:r><w:t xml:space="preserve"> = 5.413; SE = 0.746;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">t</w:t></w:r><w:r><w:t xml:space="preserve"> = 7.253;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">p</w:t></w:r><w:r><w:t xml:space="preserve"> = .001;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">γ</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t xml:space="preserve">20</w:t></w:r><w:r><w:t xml:space="preserve"> = −1.374; SE = 0.396;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">t</w:t></w:r><w:r><w:t xml:space="preserve"> = −3.470;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">p</w:t></w:r><w:r><w:t xml:space="preserve"> = .001;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">γ</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t xml:space="preserve">30</w:t></w:r><w:r><w:t xml:space="preserve"> = −1.666; SE = 0.297;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">t</w:t></w:r><w:r><w:t xml:space="preserve"> = −5.614;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">p</w:t></w:r><w:r><w:t xml:space="preserve"> = .001).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/></w:pPr><w:r><w:t xml:space="preserve">(Reporting nested model comparisons)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/></w:pPr><w:r><w:t xml:space="preserve">There was (was not) a main effect of XXX on the quadratic time term</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">(</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">χ</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t xml:space="preserve">2</w:t></w:r><w:r><w:t xml:space="preserve">(1) = 0.035,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">p</w:t></w:r><w:r><w:t xml:space="preserve"> = .852).</w:t></w:r></w:p>'

$insertPoint.InsertXML($frag)

Write-Output "Inserted example write-up paragraphs. Paragraph count now:"
Write-Output $d.Paragraphs.Count
